$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, pushing existing rows 13-63 down to 14-64
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with the new weekly record, matching the constant
# columns used throughout the rest of the sheet.
$ws.Cells.Item(13, 1).Value = 1
$ws.Cells.Item(13, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(13, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(13, 4).Value = 44575
$ws.Cells.Item(13, 5).Value = 15
$ws.Cells.Item(13, 6).Value = 100112040
$ws.Cells.Item(13, 7).Value = "Cilantro"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 250
$ws.Cells.Item(13, 11).Value = 3000
$ws.Cells.Item(13, 12).Value = 3500
$ws.Cells.Item(13, 13).Value = 3250
$ws.Cells.Item(13, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(13, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(13, 16).Value = 1625
$ws.Cells.Item(13, 17).Value = 2
$ws.Cells.Item(13, 18).Value = "Hortaliza"
